$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Title cell (A1): rewording + becomes the newest shared string entry
$ws.Range("A1").Value = "Caracteristicas de Programas de Mapas Mentais (Pesquisar e selecionar)"

# 2) Remove the defined Print_Area name
$ws.PageSetup.PrintArea = ""

# 3) Row 1 height
$ws.Rows(1).RowHeight = 34.5

# 4) Column D width (~39.14 characters)
$ws.Columns("D").ColumnWidth = 38.33

# 5) View: zoom out to 80% (page-break preview) and move the selection to G1
$excel.ActiveWindow.Zoom = 80
$ws.Range("G1").Select() | Out-Null

# 6) Page setup: scale print output to 95%
$ws.PageSetup.Zoom = 95
